$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

# Row 2
Set-TextValue 2 4 "59.307.40"
$ws.Cells.Item(2, 5).Value = "  +0.69%  "

# Row 3
Set-TextValue 3 4 "2.528.63"
$ws.Cells.Item(3, 5).Value = "  +0.69%  "

# Row 4
Set-TextValue 4 4 "0.999"
$ws.Cells.Item(4, 5).Value = "  -0.14%  "

# Row 5
Set-TextValue 5 4 "540.32"
$ws.Cells.Item(5, 5).Value = "  +1.48%  "

# Row 6
Set-TextValue 6 4 "138.50"
$ws.Cells.Item(6, 5).Value = "  -0.04%  "

# Row 7
Set-TextValue 7 4 "0.999"
$ws.Cells.Item(7, 5).Value = "  -0.09%  "

# Row 8
Set-TextValue 8 4 "0.568"
$ws.Cells.Item(8, 5).Value = "  +0.99%  "

# Row 9
Set-TextValue 9 4 "2.526.71"
$ws.Cells.Item(9, 5).Value = "  +0.59%  "

# Row 10
$ws.Cells.Item(10, 5).Value = "  +1.69%  "

# Row 11
$ws.Cells.Item(11, 5).Value = "  -0.49%  "

# Row 12
Set-TextValue 12 4 "5.39"
$ws.Cells.Item(12, 5).Value = "  -0.56%  "

# Row 13
$ws.Cells.Item(13, 5).Value = "  -1.66%  "

# Row 14
Set-TextValue 14 4 "2.961.63"
$ws.Cells.Item(14, 5).Value = "  +0.12%  "

# Row 15
Set-TextValue 15 4 "23.28"
$ws.Cells.Item(15, 5).Value = "  +1.16%  "

# Row 16
Set-TextValue 16 4 "59.188.73"
$ws.Cells.Item(16, 5).Value = "  +0.60%  "

# Row 17
$ws.Cells.Item(17, 5).Value = "  +0.53%  "

# Row 18
Set-TextValue 18 4 "2.529.05"
$ws.Cells.Item(18, 5).Value = "  +0.50%  "

# Row 19
Set-TextValue 19 4 "11.16"
$ws.Cells.Item(19, 5).Value = "  +1.44%  "

# Row 20
$ws.Cells.Item(20, 5).Value = "  +1.20%  "

# Row 21
Set-TextValue 21 4 "326.98"
$ws.Cells.Item(21, 5).Value = "  +1.53%  "

# Row 22
$ws.Cells.Item(22, 5).Value = "  -0.02%  "

# Row 23
Set-TextValue 23 4 "5.99"
$ws.Cells.Item(23, 5).Value = "  +3.27%  "

# Row 24
Set-TextValue 24 4 "65.46"
$ws.Cells.Item(24, 5).Value = "  +5.46%  "

# Row 25
$ws.Cells.Item(25, 5).Value = "  +0.43%  "

# Row 26
$ws.Cells.Item(26, 5).Value = "  +0.49%  "

# Row 27
Set-TextValue 27 4 "1.00"
$ws.Cells.Item(27, 5).Value = "  +0.05%  "

# Row 28
Set-TextValue 28 4 "7.71"
$ws.Cells.Item(28, 5).Value = "  -0.24%  "

# Row 29
Set-TextValue 29 4 "0.0₃0782"
$ws.Cells.Item(29, 5).Value = "  +2.08%  "

# Row 30
$ws.Cells.Item(30, 5).Value = "  +1.73%  "

# Row 31
Set-TextValue 31 4 "1.79"
$ws.Cells.Item(31, 5).Value = "  +0.71%  "

# Row 32
Set-TextValue 32 4 "168.31"
$ws.Cells.Item(32, 5).Value = "  +2.49%  "

# Row 33
Set-TextValue 33 4 "1.19"
$ws.Cells.Item(33, 5).Value = "  +7.69%  "

# Row 34
$ws.Cells.Item(34, 2).Value = "USDe"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue 34 4 "0.999"
$ws.Cells.Item(34, 5).Value = "  +0.00%  "

# Row 35
$ws.Cells.Item(35, 2).Value = "ImmutableX"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue 35 4 "1.47"
$ws.Cells.Item(35, 5).Value = "  +3.40%  "

# Row 36
Set-TextValue 36 4 "18.53"
$ws.Cells.Item(36, 5).Value = "  +0.62%  "

# Row 37
$ws.Cells.Item(37, 5).Value = "  -1.35%  "

# Row 38
$ws.Cells.Item(38, 5).Value = "  +0.80%  "

# Row 39
Set-TextValue 39 4 "36.82"
$ws.Cells.Item(39, 5).Value = "  -0.10%  "

# Row 40
Set-TextValue 40 4 "0.825"
$ws.Cells.Item(40, 5).Value = "  +3.30%  "

# Row 41
Set-TextValue 41 4 "3.66"
$ws.Cells.Item(41, 5).Value = "  +0.80%  "

# Row 42
Set-TextValue 42 4 "285.42"
$ws.Cells.Item(42, 5).Value = "  +2.57%  "

# Row 43
Set-TextValue 43 4 "5.24"
$ws.Cells.Item(43, 5).Value = "  +1.46%  "

# Row 44
$ws.Cells.Item(44, 5).Value = "  -0.03%  "

# Row 45
Set-TextValue 45 4 "131.88"
$ws.Cells.Item(45, 5).Value = "  +8.05%  "

# Row 46
$ws.Cells.Item(46, 5).Value = "  +2.17%  "

# Row 47
$ws.Cells.Item(47, 5).Value = "  +0.18%  "

# Row 48
$ws.Cells.Item(48, 5).Value = "  +0.25%  "

# Row 49
$ws.Cells.Item(49, 5).Value = "  +0.53%  "

# Row 50
Set-TextValue 50 4 "0.0222"
$ws.Cells.Item(50, 5).Value = "  +0.08%  "

# Row 51
Set-TextValue 51 4 "17.51"
$ws.Cells.Item(51, 5).Value = "  -0.30%  "
